# Update stock data in the workbook (data/stock.xlsx)
#
# 1) Rename the product in row 96 from "Zoro Wano" to "Zoro [Wano]"
#    (its image name "Zoro Wano.jpg" is left unchanged).
# 2) Decrement the "Stock" (column B) count by 1 for a number of rows
#    (reflects units sold / shrinkage in inventory).
# 3) Update the worksheet scroll position / active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1) Rename product -----------------------------------------------
$ws.Range("A96").Value = "Zoro [Wano]"

# --- 2) Stock count decrements ----------------------------------------
$stockUpdates = @{
    40  = 2
    68  = 1
    96  = 1
    108 = 1
    116 = 0
    118 = 0
    124 = 0
    132 = 1
    148 = 2
    153 = 0
    186 = 0
    194 = 1
    196 = 0
    255 = 4
    264 = 2
    266 = 2
    273 = 0
    274 = 2
    275 = 2
}

foreach ($row in $stockUpdates.Keys) {
    $ws.Cells.Item($row, 2).Value = $stockUpdates[$row]
}

# --- 3) Update view / selection ---------------------------------------
[void]$ws.Activate()
[void]$ws.Range("B256").Select()

$win = $excel.ActiveWindow
$win.ScrollRow = 250
$win.ScrollColumn = 1
